# CM-48.docx edit: "modificados los story points de CM-21 y CM-48 a 3 ambas"
#
# 1) Collapse the three user-story / acceptance-criteria paragraphs that
#    were split across multiple runs (with w:proofErr spell/grammar-check
#    markers in between) back down to a single run each. Replacing the
#    whole paragraph text via Find/Execute makes Word rewrite the
#    paragraph as a single run, which naturally drops the now-orphaned
#    w:proofErr markers that used to sit between the old runs.
# 2) Change the story-points value in the summary table from 5 to 3.

$d = $word.ActiveDocument

# Curly double quotation marks (U+201C / U+201D) used around "@" in the
# acceptance criteria text.
$openQuote = [char]0x201C
$closeQuote = [char]0x201D

# --- Paragraph: user story wording ---
# "...redacción " + "del mismo" + "." -> single run "...redacción del mismo."
$text1 = "Como ESTUDIANTE quiero mencionar a otro usuario en un post propio PARA hacer referencia a otro estudiante en la redacción del mismo."
$d.Content.Find.Execute($text1, $true, $false, $false, $false, $false, $true, 1, $false, $text1, 2) | Out-Null

# --- Acceptance criterion 1 ---
$text2 = "DEBE ingresar el caracter " + $openQuote + "@" + $closeQuote + " para hacer mención a otro usuario."
$d.Content.Find.Execute($text2, $true, $false, $false, $false, $false, $true, 1, $false, $text2, 2) | Out-Null

# --- Acceptance criterion 2 ---
$text3 = "Una vez ingresado el caracter " + $openQuote + "@" + $closeQuote + " seguido de tres caracteres alfanuméricos se DEBE mostrar recomendaciones de usuarios."
$d.Content.Find.Execute($text3, $true, $false, $false, $false, $false, $true, 1, $false, $text3, 2) | Out-Null

# --- Story points cell (row 1, column 2 of the summary table): 5 -> 3 ---
$table = $d.Tables.Item(1)
$pointsCell = $table.Cell(1, 2)
$pointsCell.Range.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2) | Out-Null
